$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the pre-filled sequential numbers (2..50) that used to live in B6:B54,
# leaving the cells empty but keeping their existing style (s="4").
$ws.Range("B6:B54").Value = ""

# Reflect the new selection left behind in the sheet view: B6:B54 selected,
# with B6 as the active cell (previously C11 was selected).
$ws.Range("B6:B54").Select()
